$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New repeating 8-row cycle of names for column B (PC) and column C (Mikser),
# applied to rows 2..53.
$patternB = @('Henrik','Marco','Patrick','Ledig','Martin','Louis','Marcus','Silas')
$patternC = @('Martin','Louis','Marcus','Silas','Henrik','Marco,','Patrick','Ledig')

for ($r = 2; $r -le 53; $r++) {
    $idx = ($r - 2) % 8
    $ws.Cells.Item($r, 2).Value = $patternB[$idx]
    $ws.Cells.Item($r, 3).Value = $patternC[$idx]
}
